# "fixed bug and added upload": append the newly-uploaded attendance record
# as a new row at the bottom of the 2020-11-10 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020-11-10")

$lastRow = 16
$newRow = $lastRow + 1

# Carry the existing row's formatting (bold/centered/bordered Sr. No column,
# plain data columns) down onto the new row before filling in the values.
$ws.Range("A" + $lastRow + ":G" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":G" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 1
$ws.Cells.Item($newRow, 2).Value = "sachin"
$ws.Cells.Item($newRow, 3).Value = "xyz/xyz"
$ws.Cells.Item($newRow, 4).Value = "coder"
$ws.Cells.Item($newRow, 5).Value = "21:17:22"
$ws.Cells.Item($newRow, 6).Value = 96.5882061718065
$ws.Cells.Item($newRow, 7).Value = 58.96413616731667
